# Apply the change described by the diff:
# A new observation row is inserted at row 77 (pushing the existing
# rows 77-114 down to rows 78-115). The new row duplicates the data
# that was previously in row 77, except for the date (column D) and
# the volume (column J), which get new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 77; this shifts rows
# 77:114 down to 78:115, the same way Excel's native "insert row"
# command does (including carrying the formatting/number format of
# the row being pushed down, e.g. the date format on column D).
$ws.Rows.Item(77).Insert()

# The row that used to be at 77 is now at row 78; copy its values
# into the freshly inserted (and currently empty) row 77 so that
# every column that stays constant across the dataset is carried
# over correctly, then overwrite the two columns that actually
# differ for the new observation.
$srcRow = 78
$dstRow = 77

for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item($dstRow, $col).Value2 = $ws.Cells.Item($srcRow, $col).Value2
}

# New row-specific values (column D = Fecha, column J = Volumen)
$ws.Cells.Item($dstRow, 4).Value2 = 44455
$ws.Cells.Item($dstRow, 10).Value2 = 300
